$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.252699999999995
$ws.Range("D6").Value = -7.620899999999999
$ws.Range("D7").Value = -7.062699999999998
$ws.Range("B8").Value = 5.848699999999999
$ws.Range("D8").Value = -8.163900000000005
$ws.Range("A12").Value = -22.63850000000001
$ws.Range("B12").Value = 5.2535
$ws.Range("B14").Value = 8.833600000000002
$ws.Range("D19").Value = -8.68099999999999
$ws.Range("D21").Value = -7.693799999999998
$ws.Range("B22").Value = 5.240100000000004
$ws.Range("D24").Value = -8.126799999999992
